$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + '26.459.34'
$ws.Range("E2").Value = "'" + '  -8.96%  '
$ws.Range("D3").Value = "'" + '1.621.89'
$ws.Range("E3").Value = "'" + '  -10.76%  '
$ws.Range("D4").Value = "'" + '1.043'
$ws.Range("E4").Value = "'" + '  +3.72%  '
$ws.Range("D5").Value = "'" + '220.89'
$ws.Range("E5").Value = "'" + '  -5.10%  '
$ws.Range("B6").Value = 'USDC'
$ws.Range("C6").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D6").Value = "'" + '1.038'
$ws.Range("E6").Value = "'" + '  +3.14%  '
$ws.Range("B7").Value = 'XRP'
$ws.Range("C7").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D7").Value = "'" + '0.4884'
$ws.Range("E7").Value = "'" + '  -16.62%  '
$ws.Range("D8").Value = "'" + '0.2506'
$ws.Range("E8").Value = "'" + '  -7.97%  '
$ws.Range("D9").Value = "'" + '20.78'
$ws.Range("E9").Value = "'" + '  -9.23%  '
$ws.Range("D10").Value = "'" + '0.05844'
$ws.Range("E10").Value = "'" + '  -13.59%  '
$ws.Range("D11").Value = "'" + '0.07420'
$ws.Range("E11").Value = "'" + '  -1.25%  '
$ws.Range("D12").Value = "'" + '1.720.86'
$ws.Range("E12").Value = "'" + '  -5.11%  '
$ws.Range("D13").Value = "'" + '4.462'
$ws.Range("E13").Value = "'" + '  -3.76%  '
$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D14").Value = "'" + '1.945.48'
$ws.Range("E14").Value = "'" + '  -5.69%  '
$ws.Range("B15").Value = 'Polygon'
$ws.Range("C15").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D15").Value = "'" + '0.5534'
$ws.Range("E15").Value = "'" + '  -10.71%  '
$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").Value = "'" + '0.000007728'
$ws.Range("E16").Value = "'" + '  -15.29%  '
$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").Value = "'" + '27.026.24'
$ws.Range("E17").Value = "'" + '  -6.28%  '
$ws.Range("B18").Value = 'Litecoin'
$ws.Range("C18").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D18").Value = "'" + '60.24'
$ws.Range("E18").Value = "'" + '  -19.10%  '
$ws.Range("B19").Value = 'Dai'
$ws.Range("C19").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D19").Value = "'" + '1.007'
$ws.Range("E19").Value = "'" + '  +0.06%  '
$ws.Range("B20").Value = 'Uniswap'
$ws.Range("C20").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D20").Value = "'" + '4.798'
$ws.Range("E20").Value = "'" + '  -10.99%  '
$ws.Range("D21").Value = "'" + '10.33'
$ws.Range("E21").Value = "'" + '  -8.80%  '
$ws.Range("D22").Value = "'" + '175.25'
$ws.Range("E22").Value = "'" + '  -15.19%  '
$ws.Range("B23").Value = 'BinanceUSD'
$ws.Range("C23").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D23").Value = "'" + '1.011'
$ws.Range("E23").Value = "'" + '  +0.12%  '
$ws.Range("B24").Value = 'Chainlink'
$ws.Range("C24").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D24").Value = "'" + '5.998'
$ws.Range("E24").Value = "'" + '  -10.91%  '
$ws.Range("D25").Value = "'" + '144.03'
$ws.Range("E25").Value = "'" + '  -6.39%  '
$ws.Range("D26").Value = "'" + '7.418'
$ws.Range("E26").Value = "'" + '  -4.19%  '
$ws.Range("D27").Value = "'" + '0.1112'
$ws.Range("E27").Value = "'" + '  -11.53%  '
$ws.Range("D28").Value = "'" + '14.38'
$ws.Range("E28").Value = "'" + '  -11.03%  '
$ws.Range("D29").Value = "'" + '1.358'
$ws.Range("E29").Value = "'" + '  -3.68%  '
$ws.Range("D30").Value = "'" + '0.05559'
$ws.Range("E30").Value = "'" + '  -11.76%  '
$ws.Range("D31").Value = "'" + '1.364'
$ws.Range("E31").Value = "'" + '  -4.35%  '
$ws.Range("B32").Value = 'InternetComputer(DFINITY)'
$ws.Range("C32").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D32").Value = "'" + '3.340'
$ws.Range("E32").Value = "'" + '  -8.48%  '
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").Value = "'" + '3.187'
$ws.Range("E33").Value = "'" + '  -13.29%  '
$ws.Range("D34").Value = "'" + '1.522'
$ws.Range("E34").Value = "'" + '  -9.32%  '
$ws.Range("B35").Value = 'HuobiToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D35").Value = "'" + '2.456'
$ws.Range("E35").Value = "'" + '  -3.20%  '
$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D36").Value = "'" + '0.9478'
$ws.Range("E36").Value = "'" + '  -8.99%  '
$ws.Range("B37").Value = 'MXToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D37").Value = "'" + '2.664'
$ws.Range("E37").Value = "'" + '  -2.88%  '
$ws.Range("B38").Value = 'ImmutableX'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D38").Value = "'" + '0.5633'
$ws.Range("E38").Value = "'" + '  -11.02%  '
$ws.Range("D39").Value = "'" + '0.01580'
$ws.Range("E39").Value = "'" + '  -6.85%  '
$ws.Range("B40").Value = 'TrustWalletToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D40").Value = "'" + '0.8717'
$ws.Range("E40").Value = "'" + '  +1.08%  '
$ws.Range("B41").Value = 'PaxDollar'
$ws.Range("C41").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D41").Value = "'" + '1.017'
$ws.Range("E41").Value = "'" + '  +1.07%  '
$ws.Range("B42").Value = 'Maker'
$ws.Range("C42").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D42").Value = "'" + '1.032.94'
$ws.Range("E42").Value = "'" + '  -8.20%  '
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").Value = "'" + '5.682'
$ws.Range("E43").Value = "'" + '  -10.72%  '
$ws.Range("D44").Value = "'" + '96.20'
$ws.Range("E44").Value = "'" + '  -3.38%  '
$ws.Range("D45").Value = "'" + '1.865.23'
$ws.Range("E45").Value = "'" + '  -5.39%  '
$ws.Range("B46").Value = 'BabyDogeCoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D46").Value = "'" + '0.00000000112'
$ws.Range("E46").Value = "'" + '  -0.74%  '
$ws.Range("B47").Value = 'Mantle'
$ws.Range("C47").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D47").Value = "'" + '0.4467'
$ws.Range("E47").Value = "'" + '  -1.37%  '
$ws.Range("D48").Value = "'" + '0.05316'
$ws.Range("E48").Value = "'" + '  -3.05%  '
$ws.Range("B49").Value = 'Frax'
$ws.Range("C49").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D49").Value = "'" + '0.9894'
$ws.Range("E49").Value = "'" + '  -2.31%  '
$ws.Range("D50").Value = "'" + '7.939'
$ws.Range("E50").Value = "'" + '  -2.82%  '
$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D51").Value = "'" + '53.16'
$ws.Range("E51").Value = "'" + '  -11.07%  '
